# Updates cryptos list values (price + 1h volume/change columns, and a
# 3-way row rotation for WEMIXToken / TrustWalletToken / MXToken) to match
# the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = [ordered]@{
    "D2" = "26.224.19"
    "E2" = "  -0.51%  "
    "D3" = "1.589.41"
    "E3" = "  -0.24%  "
    "D4" = "1.00"
    "E4" = "  -0.05%  "
    "D5" = "211.99"
    "D6" = "0.502"
    "E6" = "  -0.54%  "
    "D7" = "1.00"
    "E7" = "  -0.04%  "
    "D8" = "0.245"
    "E8" = "  -0.16%  "
    "E9" = "  -1.17%  "
    "D10" = "19.22"
    "E10" = "  -2.17%  "
    "D11" = "0.0847"
    "E11" = "  +0.34%  "
    "D12" = "1.813.73"
    "E12" = "  -0.13%  "
    "D13" = "1.599.22"
    "E13" = "  +0.76%  "
    "D14" = "4.01"
    "E14" = "  -1.53%  "
    "D15" = "0.517"
    "E15" = "  -0.35%  "
    "D16" = "64.05"
    "E16" = "  -0.99%  "
    "D17" = "26.227.29"
    "E17" = "  -0.51%  "
    "D18" = "0.0₃0725"
    "E18" = "  -0.67%  "
    "D19" = "214.60"
    "E19" = "  +1.15%  "
    "D20" = "7.30"
    "E20" = "  -2.55%  "
    "D21" = "0.999"
    "E21" = "  -0.16%  "
    "E22" = "  -1.03%  "
    "D23" = "2.18"
    "E23" = "  -0.71%  "
    "D24" = "8.96"
    "E24" = "  +0.26%  "
    "D25" = "144.02"
    "E25" = "  -0.81%  "
    "E26" = "  -0.06%  "
    "E27" = "  -0.91%  "
    "E28" = "  -0.95%  "
    "D29" = "15.14"
    "E29" = "  -1.01%  "
    "D30" = "0.0496"
    "E30" = "  -1.87%  "
    "E31" = "  +0.63%  "
    "D32" = "3.19"
    "E32" = "  -1.31%  "
    "D33" = "1.389.17"
    "E33" = "  +6.71%  "
    "D34" = "2.93"
    "E34" = "  -1.87%  "
    "E35" = "  -0.19%  "
    "E36" = "  -1.46%  "
    "D37" = "0.583"
    "E37" = "  -5.58%  "
    "E38" = "  -0.76%  "
    "D39" = "0.820"
    "E39" = "  +0.39%  "
    "D40" = "5.85"
    "E40" = "  +4.05%  "
    "D41" = "0.999"
    "E41" = "  -0.12%  "
    "B42" = "TrustWalletToken"
    "C42" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
    "D42" = "0.769"
    "E42" = "  +0.84%  "
    "B43" = "MXToken"
    "C43" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D43" = "2.13"
    "E43" = "  -0.19%  "
    "B44" = "WEMIXToken"
    "C44" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D44" = "0.915"
    "E44" = "  -17.30%  "
    "D45" = "1.725.75"
    "E45" = "  -0.16%  "
    "D46" = "61.08"
    "E46" = "  -2.84%  "
    "D47" = "86.10"
    "E47" = "  -2.83%  "
    "E48" = "  -1.34%  "
    "D49" = "0.0979"
    "E49" = "  -1.46%  "
    "D50" = "0.0500"
    "E50" = "  -1.07%  "
    "D51" = "0.997"
    "E51" = "  -0.26%  "
}

foreach ($cell in $edits.Keys) {
    # Leading apostrophe forces text interpretation so values such as
    # "1.00" or "26.224.19" are not coerced into numbers and lose their
    # formatting (mirrors how these sheets are generated as text cells).
    $ws.Range($cell).Value = "'" + $edits[$cell]
    $ws.Range($cell).Style = "Normal"
}
